$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data.
# Values are written with a leading apostrophe to force text
# interpretation (prices like "243.41" would otherwise be parsed
# as numbers by Excel), then ClearFormats() strips the resulting
# quote-prefix style so no stray style index is introduced.
$updates = @(
    @{Cell="D2"; Value="29.434.57"}
    @{Cell="E2"; Value="  +0.22%  "}
    @{Cell="D3"; Value="1.867.46"}
    @{Cell="E3"; Value="  -0.78%  "}
    @{Cell="E4"; Value="  -0.01%  "}
    @{Cell="D5"; Value="243.41"}
    @{Cell="E5"; Value="  +0.15%  "}
    @{Cell="D6"; Value="0.7039"}
    @{Cell="E6"; Value="  -1.14%  "}
    @{Cell="E7"; Value="  -0.02%  "}
    @{Cell="D8"; Value="0.3136"}
    @{Cell="E8"; Value="  -0.79%  "}
    @{Cell="D9"; Value="0.07845"}
    @{Cell="E9"; Value="  -2.09%  "}
    @{Cell="D10"; Value="24.45"}
    @{Cell="E10"; Value="  -2.63%  "}
    @{Cell="D11"; Value="0.08021"}
    @{Cell="E11"; Value="  -3.72%  "}
    @{Cell="D12"; Value="1.894.65"}
    @{Cell="E12"; Value="  +0.27%  "}
    @{Cell="D13"; Value="5.190"}
    @{Cell="E13"; Value="  -1.59%  "}
    @{Cell="D14"; Value="93.29"}
    @{Cell="E14"; Value="  -1.69%  "}
    @{Cell="D15"; Value="0.6999"}
    @{Cell="E15"; Value="  -2.62%  "}
    @{Cell="D16"; Value="6.452"}
    @{Cell="E16"; Value="  +1.30%  "}
    @{Cell="D17"; Value="29.501.56"}
    @{Cell="E17"; Value="  +0.39%  "}
    @{Cell="D18"; Value="0.000008311"}
    @{Cell="E18"; Value="  -4.13%  "}
    @{Cell="D19"; Value="254.06"}
    @{Cell="E19"; Value="  +4.52%  "}
    @{Cell="D20"; Value="2.141.34"}
    @{Cell="E20"; Value="  -0.59%  "}
    @{Cell="D21"; Value="13.12"}
    @{Cell="E21"; Value="  -1.69%  "}
    @{Cell="E22"; Value="  -0.07%  "}
    @{Cell="D23"; Value="7.594"}
    @{Cell="E23"; Value="  -3.29%  "}
    @{Cell="D24"; Value="1.002"}
    @{Cell="E24"; Value="  -0.06%  "}
    @{Cell="D25"; Value="0.1552"}
    @{Cell="E25"; Value="  -1.42%  "}
    @{Cell="D26"; Value="9.025"}
    @{Cell="E26"; Value="  -0.82%  "}
    @{Cell="D27"; Value="160.75"}
    @{Cell="E27"; Value="  -1.54%  "}
    @{Cell="D28"; Value="18.76"}
    @{Cell="E28"; Value="  +0.77%  "}
    @{Cell="D29"; Value="1.500"}
    @{Cell="E29"; Value="  -0.81%  "}
    @{Cell="D30"; Value="4.320"}
    @{Cell="E30"; Value="  -2.84%  "}
    @{Cell="D31"; Value="4.253"}
    @{Cell="E31"; Value="  -2.33%  "}
    @{Cell="D32"; Value="1.201"}
    @{Cell="E32"; Value="  -0.49%  "}
    @{Cell="D33"; Value="0.05299"}
    @{Cell="E33"; Value="  -1.90%  "}
    @{Cell="D34"; Value="1.884"}
    @{Cell="E34"; Value="  -3.32%  "}
    @{Cell="D35"; Value="0.7445"}
    @{Cell="E35"; Value="  -3.98%  "}
    @{Cell="E36"; Value="  -1.93%  "}
    @{Cell="D37"; Value="2.722"}
    @{Cell="E37"; Value="  +1.51%  "}
    @{Cell="D38"; Value="0.01869"}
    @{Cell="E38"; Value="  -1.28%  "}
    @{Cell="D39"; Value="1.258.45"}
    @{Cell="E39"; Value="  -1.23%  "}
    @{Cell="D40"; Value="2.742"}
    @{Cell="E40"; Value="  -0.23%  "}
    @{Cell="D41"; Value="0.8984"}
    @{Cell="E41"; Value="  -2.26%  "}
    @{Cell="D42"; Value="108.65"}
    @{Cell="E42"; Value="  -4.01%  "}
    @{Cell="D43"; Value="5.954"}
    @{Cell="E43"; Value="  -8.78%  "}
    @{Cell="D44"; Value="71.19"}
    @{Cell="E44"; Value="  -4.54%  "}
    @{Cell="E45"; Value="  -0.02%  "}
    @{Cell="D46"; Value="0.00000000129"}
    @{Cell="E46"; Value="  +0.35%  "}
    @{Cell="D47"; Value="2.039.38"}
    @{Cell="E47"; Value="  -0.51%  "}
    @{Cell="B48"; Value="Mantle"}
    @{Cell="C48"; Value="https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"}
    @{Cell="D48"; Value="0.5191"}
    @{Cell="E48"; Value="  -0.62%  "}
    @{Cell="B49"; Value="RenderToken"}
    @{Cell="C49"; Value="https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"}
    @{Cell="D49"; Value="1.792"}
    @{Cell="E49"; Value="  -1.40%  "}
    @{Cell="D50"; Value="9.494"}
    @{Cell="E50"; Value="  -0.89%  "}
    @{Cell="D51"; Value="0.4301"}
    @{Cell="E51"; Value="  -1.90%  "}
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    $rng.Value = "'" + $u.Value
    $rng.ClearFormats()
}
